# Mark the "add your picture and name to the team section" row (OmarAlcaesar)
# as Done in the status column, matching the commit:
# "adding my picture and name in to the team section"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "Done"

# Reflect the cursor/selection position left by the edit.
$ws.Range("C18").Select() | Out-Null
